$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert three new columns before column F (shifting old F -> I)
$ws.Range("F:H").Insert()

# New header/data values (order matches original authoring sequence)
$ws.Range("F1").Value = "Product Category"
$ws.Range("F3").Value = "Books & Media"
$ws.Range("G1").Value = "Product Subcategory"
$ws.Range("H1").Value = "Exact Product Name"
$ws.Range("G3").Value = "Literature & Fiction"
$ws.Range("H3").Value = "Shri Ramcharitmanas"

# Set column widths
$ws.Range("F1").ColumnWidth = 18.28515625
$ws.Range("G1").ColumnWidth = 20.42578125
$ws.Range("H1").ColumnWidth = 19.140625

# Apply the quote-prefixed bordered style (matching E3) to the new G/H data cells
$ws.Range("E3").Copy()
$ws.Range("G2:H10").PasteSpecial(-4122)

# Selection / view
$ws.Range("H3").Select()
